# Updates the crypto price/volume snapshot (and two rank swaps: rows 19/20
# and 27/28) to match the refreshed data pulled by the scheduled scraper.
# Cells that would otherwise be auto-coerced to a number by the smart-entry
# parser (single '.' decimal separators) are forced back to text via
# NumberFormat "@" before the value is written, matching the original
# inlineStr/text cell type.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.840.19"
$ws.Range("E2").Value = "  +1.12%  "
$ws.Range("D3").Value = "2.623.80"
$ws.Range("E3").Value = "  +1.13%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "598.64"
$ws.Range("E5").Value = "  +0.63%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.88"
$ws.Range("E6").Value = "  -0.33%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.547"
$ws.Range("E8").Value = "  +0.62%  "
$ws.Range("D9").Value = "2.622.54"
$ws.Range("E9").Value = "  +1.09%  "
$ws.Range("E10").Value = "  +8.57%  "
$ws.Range("E11").Value = "  +0.91%  "
$ws.Range("E12").Value = "  +0.61%  "
$ws.Range("E13").Value = "  -1.57%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.71"
$ws.Range("E14").Value = "  -2.90%  "
$ws.Range("E15").Value = "  +2.78%  "
$ws.Range("D16").Value = "3.102.00"
$ws.Range("E16").Value = "  +1.54%  "
$ws.Range("D17").Value = "67.743.23"
$ws.Range("E17").Value = "  +1.01%  "
$ws.Range("D18").Value = "2.629.45"
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.19"
$ws.Range("E19").Value = "  -1.40%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "366.56"
$ws.Range("E20").Value = "  +2.53%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.63"
$ws.Range("E21").Value = "  -2.77%  "
$ws.Range("E22").Value = "  -0.73%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.04"
$ws.Range("E23").Value = "  -2.25%  "
$ws.Range("E24").Value = "  -0.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.85"
$ws.Range("E25").Value = "  -6.88%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "67.35"
$ws.Range("E26").Value = "  +0.27%  "
$ws.Range("B27").Value = "WrappedeETH"
$ws.Range("C27").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D27").Value = "2.747.29"
$ws.Range("E27").Value = "  +1.28%  "
$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000104"
$ws.Range("E28").Value = "  -0.38%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "578.31"
$ws.Range("E29").Value = "  -4.79%  "
$ws.Range("E30").Value = "  +0.51%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.42"
$ws.Range("E31").Value = "  -3.47%  "
$ws.Range("E32").Value = "  -2.56%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.86"
$ws.Range("E33").Value = "  -0.26%  "
$ws.Range("E34").Value = "  -1.89%  "
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("E36").Value = "  -3.83%  "
$ws.Range("E37").Value = "  -2.28%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "158.07"
$ws.Range("E38").Value = "  +2.24%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.369"
$ws.Range("E40").Value = "  -0.56%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.85"
$ws.Range("E41").Value = "  +1.79%  "
$ws.Range("E42").Value = "  -3.25%  "
$ws.Range("E43").Value = "  -4.13%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.25"
$ws.Range("E44").Value = "  -0.76%  "
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.42"
$ws.Range("E46").Value = "  -0.10%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "157.00"
$ws.Range("E47").Value = "  +0.26%  "
$ws.Range("D48").Value = "0.0₆0285"
$ws.Range("E48").Value = "  -7.86%  "
$ws.Range("E49").Value = "  -0.74%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "21.05"
$ws.Range("E50").Value = "  -2.42%  "
$ws.Range("E51").Value = "  +1.36%  "
